$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 370/371, pushing the existing rows 370-386 down to 372-388.
$ws.Rows("370:371").Insert()

# ---- New row 370 ----
$ws.Cells.Item(370, 1).Value = 4
$ws.Cells.Item(370, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(370, 3).Value = "Los Lagos"
$ws.Cells.Item(370, 4).Value = 45041
$ws.Cells.Item(370, 5).Value = 10
$ws.Cells.Item(370, 6).Value = 100112024
$ws.Cells.Item(370, 7).Value = "Choclo"
$ws.Cells.Item(370, 8).Value = "Choclero"
$ws.Cells.Item(370, 9).Value = "Primera"
$ws.Cells.Item(370, 10).Value = 6000
$ws.Cells.Item(370, 11).Value = 600
$ws.Cells.Item(370, 12).Value = 650
$ws.Cells.Item(370, 13).Value = 625
$ws.Cells.Item(370, 14).Value = "$/unidad"
$ws.Cells.Item(370, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(370, 16).Value = 625
$ws.Cells.Item(370, 17).Value = 1
$ws.Cells.Item(370, 18).Value = "Hortaliza"

# ---- New row 371 ----
$ws.Cells.Item(371, 1).Value = 4
$ws.Cells.Item(371, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(371, 3).Value = "Los Lagos"
$ws.Cells.Item(371, 4).Value = 45041
$ws.Cells.Item(371, 5).Value = 10
$ws.Cells.Item(371, 6).Value = 100112024
$ws.Cells.Item(371, 7).Value = "Choclo"
$ws.Cells.Item(371, 8).Value = "Dulce o Americano"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 100
$ws.Cells.Item(371, 11).Value = 16000
$ws.Cells.Item(371, 12).Value = 16000
$ws.Cells.Item(371, 13).Value = 16000
$ws.Cells.Item(371, 14).Value = "$/malla 60 unidades"
$ws.Cells.Item(371, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(371, 16).Value = 267
$ws.Cells.Item(371, 17).Value = 60
$ws.Cells.Item(371, 18).Value = "Hortaliza"
